$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.593300805548547; C = 0.798316168795588; D = 0.809737487250624; E = 0.801269845163997; F = 0.802759474182296; G = 0.803503442074255 }
    3  = @{ B = 0.891425446185947; C = 0.734343884976076; D = 0.66011845184374;  E = 0.773993520808291; F = 0.769785232844933; G = 0.734188607563879 }
    4  = @{ B = 0.837842434332639; C = 0.649556135900333; D = 0.544523480585528; E = 0.727832934920089; F = 0.722448306363765; G = 0.676610115099755 }
    5  = @{ B = 0.747059011250423; C = 0.760748646394783; D = 0.820640960368556; E = 0.745104631264865; F = 0.72988960058406;  G = 0.728092737150705 }
    6  = @{ B = 0.620917729383005; C = 0.78913119672879;  D = 0.855042207071307; E = 0.753586689991038; F = 0.834664561377067; G = 0.754698076525027 }
    7  = @{ B = 0.701622795907415; C = 0.712710646471761; D = 0.757480297774307; E = 0.637897645695564; F = 0.877803125109165; G = 0.700890373546844 }
    8  = @{ B = 0.71036942310177;  C = 0.767459227638374; D = 0.769430442315501; E = 0.758346955137085; F = 0.84096493071521;  G = 0.746768125063093 }
    9  = @{ B = 0.776451040625449; C = 0.752356847825706; D = 0.829686492804909; E = 0.705157484400394; F = 0.849273919987832; G = 0.66342378640839  }
    10 = @{ B = 0.869708886859866; C = 0.815652162844027; D = 0.834311730887263; E = 0.8156709550428;   F = 0.867891452365697; G = 0.7671375817381   }
    11 = @{ B = 0.849616378214432; C = 0.734028088236682; D = 0.756084447572338; E = 0.685545279841557; F = 0.842706846604206; G = 0.69712559445501  }
    12 = @{ B = 0.891121217137727; C = 0.819983433106146; D = 0.87483983315479;  E = 0.673226727255466; F = 0.913082349886563; G = 0.869575896966582 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
